# Regenerate the "K" column (G) values for the save_data sheet.
# These replace the previous "Strike#" derived counts with the
# recalculated K counts (std/mean based s_vals calculation).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 2
    3  = 0
    4  = 0
    5  = 2
    6  = 2
    8  = 0
    9  = 3
    10 = 0
    11 = 1
    12 = 1
    13 = 0
    14 = 1
    15 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
